$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: id-div-container (id column D) / container de tudo
$ws.Range("B4").Value = "id-div-container"
$ws.Range("D4").Value = "x"
$ws.Range("E4").Value = "container de tudo"

# Row 5: class-img-logo (class column C) / logotipo, generica
$ws.Range("B5").Value = "class-img-logo"
$ws.Range("C5").Value = "x"
$ws.Range("E5").Value = "logotipo, generica"

# Row 6: id-nav-pai-ul (id column D) / nav pai do ul e li
$ws.Range("E6").Value = "nav pai do ul e li"
$ws.Range("B6").Value = "id-nav-pai-ul"
$ws.Range("D6").Value = "x"

# Row 7: class-li (class column C) / filhos ul
$ws.Range("B7").Value = "class-li"
$ws.Range("C7").Value = "x"
$ws.Range("E7").Value = "filhos ul"

# Row 8: id-ul-pai-li (id column D) / controla a lista
$ws.Range("B8").Value = "id-ul-pai-li"
$ws.Range("D8").Value = "x"
$ws.Range("E8").Value = "controla a lista"

# Update selected cell to match the saved view state
$ws.Range("B8").Select()
